$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.635.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.628.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.32%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.655"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.614.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.07%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.665"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.47%  "
$ws.Range("E13").Value = "  +5.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.210.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.17%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.99%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.625.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.639.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.55%  "
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +5.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +16.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "619.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.118"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.93%  "
$ws.Range("E36").Value = "  +11.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.407"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.22%  "
$ws.Range("E38").Value = "  +3.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.347.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.79%  "
$ws.Range("E44").Value = "  +7.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.139"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
